$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range('D2')
$cell.NumberFormat = '@'
$cell.Value = '41.689.26'
$cell.Style = 'Normal'
$ws.Range('E2').Value = '  +0.40%  '
$cell = $ws.Range('D3')
$cell.NumberFormat = '@'
$cell.Value = '2.478.44'
$cell.Style = 'Normal'
$ws.Range('E3').Value = '  -0.03%  '
$ws.Range('E4').Value = '  -0.02%  '
$cell = $ws.Range('D5')
$cell.NumberFormat = '@'
$cell.Value = '318.99'
$cell.Style = 'Normal'
$ws.Range('E5').Value = '  +1.73%  '
$cell = $ws.Range('D6')
$cell.NumberFormat = '@'
$cell.Value = '93.12'
$cell.Style = 'Normal'
$ws.Range('E6').Value = '  -0.10%  '
$ws.Range('E7').Value = '  +2.15%  '
$cell = $ws.Range('D9')
$cell.NumberFormat = '@'
$cell.Value = '0.520'
$cell.Style = 'Normal'
$ws.Range('E9').Value = '  +2.74%  '
$cell = $ws.Range('D10')
$cell.NumberFormat = '@'
$cell.Value = '0.0885'
$cell.Style = 'Normal'
$ws.Range('E10').Value = '  +12.89%  '
$cell = $ws.Range('D11')
$cell.NumberFormat = '@'
$cell.Value = '33.11'
$cell.Style = 'Normal'
$ws.Range('E11').Value = '  +1.18%  '
$ws.Range('E12').Value = '  +0.98%  '
$cell = $ws.Range('D13')
$cell.NumberFormat = '@'
$cell.Value = '2.860.31'
$cell.Style = 'Normal'
$ws.Range('E13').Value = '  -0.16%  '
$ws.Range('E14').Value = '  +1.84%  '
$cell = $ws.Range('D15')
$cell.NumberFormat = '@'
$cell.Value = '15.71'
$cell.Style = 'Normal'
$ws.Range('E15').Value = '  -2.67%  '
$cell = $ws.Range('D16')
$cell.NumberFormat = '@'
$cell.Value = '2.485.83'
$cell.Style = 'Normal'
$ws.Range('E16').Value = '  -1.20%  '
$cell = $ws.Range('D17')
$cell.NumberFormat = '@'
$cell.Value = '0.794'
$cell.Style = 'Normal'
$ws.Range('E17').Value = '  +4.23%  '
$cell = $ws.Range('D18')
$cell.NumberFormat = '@'
$cell.Value = '41.648.40'
$cell.Style = 'Normal'
$ws.Range('E18').Value = '  +0.31%  '
$ws.Range('E19').Value = '  +2.84%  '
$cell = $ws.Range('D20')
$cell.NumberFormat = '@'
$cell.Value = '6.51'
$cell.Style = 'Normal'
$ws.Range('E20').Value = '  +1.66%  '
$cell = $ws.Range('D21')
$cell.NumberFormat = '@'
$cell.Value = '71.59'
$cell.Style = 'Normal'
$ws.Range('E21').Value = '  -0.16%  '
$cell = $ws.Range('D22')
$cell.NumberFormat = '@'
$cell.Value = '11.57'
$cell.Style = 'Normal'
$ws.Range('E22').Value = '  +2.23%  '
$cell = $ws.Range('D23')
$cell.NumberFormat = '@'
$cell.Value = '242.23'
$cell.Style = 'Normal'
$ws.Range('E23').Value = '  +2.35%  '
$ws.Range('E24').Value = '  +1.86%  '
$ws.Range('E25').Value = '  +1.93%  '
$ws.Range('E26').Value = '  -0.06%  '
$cell = $ws.Range('D27')
$cell.NumberFormat = '@'
$cell.Value = '24.93'
$cell.Style = 'Normal'
$ws.Range('E27').Value = '  +0.32%  '
$ws.Range('E28').Value = '  +3.89%  '
$cell = $ws.Range('D29')
$cell.NumberFormat = '@'
$cell.Value = '9.91'
$cell.Style = 'Normal'
$ws.Range('E29').Value = '  +2.75%  '
$cell = $ws.Range('D30')
$cell.NumberFormat = '@'
$cell.Value = '36.69'
$cell.Style = 'Normal'
$ws.Range('E30').Value = '  +2.43%  '
$cell = $ws.Range('D31')
$cell.NumberFormat = '@'
$cell.Value = '156.94'
$cell.Style = 'Normal'
$ws.Range('E31').Value = '  -0.71%  '
$cell = $ws.Range('D32')
$cell.NumberFormat = '@'
$cell.Value = '5.55'
$cell.Style = 'Normal'
$ws.Range('E32').Value = '  +1.25%  '
$ws.Range('E33').Value = '  -0.06%  '
$cell = $ws.Range('D34')
$cell.NumberFormat = '@'
$cell.Value = '0.0773'
$cell.Style = 'Normal'
$ws.Range('E34').Value = '  +2.52%  '
$ws.Range('E35').Value = '  -0.07%  '
$cell = $ws.Range('D36')
$cell.NumberFormat = '@'
$cell.Value = '17.58'
$cell.Style = 'Normal'
$ws.Range('E36').Value = '  +1.06%  '
$ws.Range('E37').Value = '  +0.57%  '
$ws.Range('E38').Value = '  +0.69%  '
$ws.Range('E39').Value = '  +1.37%  '
$ws.Range('E40').Value = '  -0.74%  '
$cell = $ws.Range('D41')
$cell.NumberFormat = '@'
$cell.Value = '4.03'
$cell.Style = 'Normal'
$ws.Range('E41').Value = '  -1.72%  '
$ws.Range('E42').Value = '  +2.10%  '
$cell = $ws.Range('D43')
$cell.NumberFormat = '@'
$cell.Value = '19.67'
$cell.Style = 'Normal'
$ws.Range('E43').Value = '  -0.30%  '
$cell = $ws.Range('D44')
$cell.NumberFormat = '@'
$cell.Value = '1.985.25'
$cell.Style = 'Normal'
$ws.Range('E44').Value = '  +0.77%  '
$ws.Range('E45').Value = '  +1.11%  '
$ws.Range('E46').Value = '  +3.30%  '
$cell = $ws.Range('D47')
$cell.NumberFormat = '@'
$cell.Value = '9.23'
$cell.Style = 'Normal'
$ws.Range('E47').Value = '  +0.95%  '
$cell = $ws.Range('D48')
$cell.NumberFormat = '@'
$cell.Value = '2.712.90'
$cell.Style = 'Normal'
$cell = $ws.Range('D49')
$cell.NumberFormat = '@'
$cell.Value = '97.76'
$cell.Style = 'Normal'
$ws.Range('E49').Value = '  -0.01%  '
$cell = $ws.Range('D50')
$cell.NumberFormat = '@'
$cell.Value = '68.24'
$cell.Style = 'Normal'
$ws.Range('E50').Value = '  +0.50%  '
$cell = $ws.Range('D51')
$cell.NumberFormat = '@'
$cell.Value = '74.22'
$cell.Style = 'Normal'
$ws.Range('E51').Value = '  +2.69%  '
